$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Shift existing rows 12-15 down to 13-16 to make room for a new "Jurisdiction" row.
# Using Copy(Destination) on the whole block (instead of Rows.Insert) avoids Excel
# fabricating a spurious extra cell style for the newly vacated row.
$ws.Range("A12:B15").Copy($ws.Range("A13:B16"))
# Copying a blank source cell does not clear a non-blank destination cell, so make
# sure B15 (which must become blank, like the old B14 "Copyright" row) is cleared.
$ws.Range("B15").ClearContents()

# Update simple value changes.
$ws.Range("B3").Value2 = "0.1.7"
$ws.Range("B6").Value2 = "draft"
$ws.Range("B8").Value2 = "2024-11-22T12:33:30-06:00"
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# New "Jurisdiction" row inserted at row 12 (empty value).
$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""

Write-Output "Done"
